# Applies the row-shuffle edit described by the diff: the dynamic columns
# (D, K, L, M, N, O, P, Q, R, S, T) of data rows 2-26 are rearranged as if
# each row's data (except the constant A/B/C/E/F/G/H/I/J columns) was moved
# to a different row. Row 7 keeps its original data (maps to itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each target row number, the list below gives the values that row
# must end up holding, taken (per the diff) from another row's original data.
$rows = @{
    2  = @{ D = 44784; K = "Hayward";         L = "Primera"; M = 300; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1083; T = 18 }
    3  = @{ D = 44291; K = "Hayward";         L = "Primera"; M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 972;  T = 18 }
    4  = @{ D = 44614; K = "Hayward";         L = "Primera"; M = 250; N = 20000; O = 21000; P = 20500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1139; T = 18 }
    5  = @{ D = 44263; K = "Hayward";         L = "Primera"; M = 250; N = 21000; O = 22000; P = 21500; Q = "`$/caja 18 kilos";           R = "Región de O'Higgins"; S = 1194; T = 18 }
    6  = @{ D = 44789; K = "Hayward";         L = "Segunda"; M = 250; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1083; T = 18 }
    7  = @{ D = 44323; K = "Hayward";         L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1194; T = 18 }
    8  = @{ D = 45086; K = "Hayward";         L = "Especial"; M = 250; N = 25000; O = 26000; P = 25500; Q = "`$/bandeja 18 kilos";       R = "Región de O'Higgins"; S = 1417; T = 18 }
    9  = @{ D = 45086; K = "Hayward";         L = "Primera"; M = 250; N = 20000; O = 21000; P = 20500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1139; T = 18 }
    10 = @{ D = 45069; K = "Sin especificar"; L = "Primera"; M = 370; N = 19000; O = 20000; P = 19486; Q = "`$/bandeja 18 kilos";        R = "Región Metropolitana"; S = 1083; T = 18 }
    11 = @{ D = 44706; K = "Hayward";         L = "Primera"; M = 400; N = 9000;  O = 10000; P = 9500;  Q = "`$/bandeja 10 kilos";        R = "Región de O'Higgins"; S = 950;  T = 10 }
    12 = @{ D = 44819; K = "Hayward";         L = "Primera"; M = 300; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 10 kilos";        R = "Región de O'Higgins"; S = 1750; T = 10 }
    13 = @{ D = 44489; K = "Hayward";         L = "Primera"; M = 300; N = 26000; O = 27000; P = 26500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1472; T = 18 }
    14 = @{ D = 44487; K = "Hayward";         L = "Primera"; M = 300; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 10 kilos";        R = "Región de O'Higgins"; S = 1450; T = 10 }
    15 = @{ D = 44307; K = "Hayward";         L = "Primera"; M = 250; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1083; T = 18 }
    16 = @{ D = 44673; K = "Hayward";         L = "Especial"; M = 400; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 10 kilos";       R = "Región de O'Higgins"; S = 1450; T = 10 }
    17 = @{ D = 44616; K = "Hayward";         L = "Segunda"; M = 300; N = 16000; O = 17000; P = 16500; Q = "`$/caja 18 kilos granel";    R = "Región de O'Higgins"; S = 917;  T = 18 }
    18 = @{ D = 44602; K = "Hayward";         L = "Primera"; M = 270; N = 20000; O = 21000; P = 20500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1139; T = 18 }
    19 = @{ D = 44491; K = "Hayward";         L = "Primera"; M = 300; N = 14000; O = 15000; P = 14500; Q = "`$/bandeja 10 kilos";        R = "Región de O'Higgins"; S = 1450; T = 10 }
    20 = @{ D = 44656; K = "Hayward";         L = "Primera"; M = 270; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1083; T = 18 }
    21 = @{ D = 44991; K = "Hayward";         L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1361; T = 18 }
    22 = @{ D = 45002; K = "Hayward";         L = "Segunda"; M = 300; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1361; T = 18 }
    23 = @{ D = 44418; K = "Hayward";         L = "Primera"; M = 240; N = 10000; O = 11000; P = 10500; Q = "`$/bandeja 10 kilos";        R = "Región de O'Higgins"; S = 1050; T = 10 }
    24 = @{ D = 44629; K = "Hayward";         L = "Segunda"; M = 300; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 972;  T = 18 }
    25 = @{ D = 45043; K = "Hayward";         L = "Segunda"; M = 300; N = 21000; O = 22000; P = 21500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1194; T = 18 }
    26 = @{ D = 45034; K = "Hayward";         L = "Primera"; M = 250; N = 25000; O = 26000; P = 25600; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins"; S = 1422; T = 18 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
